$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

# Clear G2 and G7 (remove their text content)
$ws.Range("G2").Value = ""
$ws.Range("G7").Value = ""

# Update Azami amounts in row 13
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 300 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 851,5 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

# Fill in previously empty E24 and E25
$ws.Range("E24").Value = "600 TL"
$ws.Range("E25").Value = "600 TL"
